$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-05 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-06 Friday", 2)
$d.Content.Find.Execute("313×5=", $true, $false, $false, $false, $false, $true, 1, $false, "326×5=", 2)
$d.Content.Find.Execute("131×7=", $true, $false, $false, $false, $false, $true, 1, $false, "378×9=", 2)
$d.Content.Find.Execute("609×3=", $true, $false, $false, $false, $false, $true, 1, $false, "815×3=", 2)
$d.Content.Find.Execute("636×2=", $true, $false, $false, $false, $false, $true, 1, $false, "477×2=", 2)
$d.Content.Find.Execute("533×9=", $true, $false, $false, $false, $false, $true, 1, $false, "411×2=", 2)
$d.Content.Find.Execute("439×3=", $true, $false, $false, $false, $false, $true, 1, $false, "512×9=", 2)
$d.Content.Find.Execute("631×3=", $true, $false, $false, $false, $false, $true, 1, $false, "811×2=", 2)
$d.Content.Find.Execute("483×7=", $true, $false, $false, $false, $false, $true, 1, $false, "226×5=", 2)
$d.Content.Find.Execute("560×5=", $true, $false, $false, $false, $false, $true, 1, $false, "991×3=", 2)
$d.Content.Find.Execute("235×7=", $true, $false, $false, $false, $false, $true, 1, $false, "952×4=", 2)
$d.Content.Find.Execute("742×7=", $true, $false, $false, $false, $false, $true, 1, $false, "423×2=", 2)
$d.Content.Find.Execute("813×2=", $true, $false, $false, $false, $false, $true, 1, $false, "653×3=", 2)
$d.Content.Find.Execute("400×4=", $true, $false, $false, $false, $false, $true, 1, $false, "112×3=", 2)
$d.Content.Find.Execute("404×6=", $true, $false, $false, $false, $false, $true, 1, $false, "987×5=", 2)
$d.Content.Find.Execute("484×9=", $true, $false, $false, $false, $false, $true, 1, $false, "530×2=", 2)
$d.Content.Find.Execute("730×4=", $true, $false, $false, $false, $false, $true, 1, $false, "559×9=", 2)
$d.Content.Find.Execute("428×6=", $true, $false, $false, $false, $false, $true, 1, $false, "490×2=", 2)
$d.Content.Find.Execute("791×7=", $true, $false, $false, $false, $false, $true, 1, $false, "141×8=", 2)
$d.Content.Find.Execute("641×9=", $true, $false, $false, $false, $false, $true, 1, $false, "386×2=", 2)
$d.Content.Find.Execute("289×7=", $true, $false, $false, $false, $false, $true, 1, $false, "688×8=", 2)
$d.Content.Find.Execute("921×8=", $true, $false, $false, $false, $false, $true, 1, $false, "140×3=", 2)
$d.Content.Find.Execute("413×2=", $true, $false, $false, $false, $false, $true, 1, $false, "525×4=", 2)
$d.Content.Find.Execute("401×6=", $true, $false, $false, $false, $false, $true, 1, $false, "963×2=", 2)
$d.Content.Find.Execute("735×6=", $true, $false, $false, $false, $false, $true, 1, $false, "530×9=", 2)
$d.Content.Find.Execute("167×7=", $true, $false, $false, $false, $false, $true, 1, $false, "102×4=", 2)
